$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (Periodo Mora) on rows 16-28: the list of periods is reordered
# from ascending (2001..2012,2101) to descending (2101,2012..2001).
$ws.Cells.Item(16, 5).Value2 = "2101"
$ws.Cells.Item(17, 5).Value2 = "2012"
$ws.Cells.Item(18, 5).Value2 = "2011"
$ws.Cells.Item(19, 5).Value2 = "2010"
$ws.Cells.Item(20, 5).Value2 = "2009"
$ws.Cells.Item(21, 5).Value2 = "2008"
$ws.Cells.Item(22, 5).Value2 = "2007"
$ws.Cells.Item(23, 5).Value2 = "2006"
$ws.Cells.Item(24, 5).Value2 = "2005"
$ws.Cells.Item(25, 5).Value2 = "2004"
$ws.Cells.Item(26, 5).Value2 = "2003"
$ws.Cells.Item(27, 5).Value2 = "2002"
$ws.Cells.Item(28, 5).Value2 = "2001"

# Column F (Valor Mora) values for the rows that swapped periods (row 16 now
# holds what used to be row 28's period, and vice versa), so their amounts
# follow the period.
$ws.Cells.Item(16, 6).Value2 = 26500
$ws.Cells.Item(28, 6).Value2 = 33125
